$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "69.746.26"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.00%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.531.39"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.08%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "604.64"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "193.97"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  -0.58%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -4.83%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.646"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.88%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "53.14"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.61%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0000302"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.41%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "9.47"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.40%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.092.90"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.89%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "591.00"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.78%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "12.76"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.82%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "69.862.56"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "18.96"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.530.27"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("E20").Value = "  +1.75%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.980"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "17.71"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.72%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "103.13"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("E24").Value = "  +0.81%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "4.62"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.58%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.03"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.77%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.69"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.36%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.49"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.06%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "33.07"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.80%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.99"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.22"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.14%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "12.28"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.83%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.114"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.64%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "63.25"
$c.Style = "Normal"
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +5.69%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.815.06"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.14%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0₃0816"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +3.21%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "511.25"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  -0.69%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "36.38"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("E43").Value = "  -2.37%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0445"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.29%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.33"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("E48").Value = "  +0.07%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.47"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("E51").Value = "  +1.37%  "
